$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values are plain numerics,
# so the engine keeps them as exact text instead of coercing to floating point.
$textCells = @("D5", "D6", "D10", "D11", "D14", "D15", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D41", "D42", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.820.89'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.478.24'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '318.98'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").Value = '93.36'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("D10").Value = '0.0883'
$ws.Range("E10").Value = '  +11.64%  '
$ws.Range("D11").Value = '33.34'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '2.861.33'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = '6.96'
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("D15").Value = '15.68'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '2.479.76'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("E17").Value = '  +3.41%  '
$ws.Range("D18").Value = '41.787.50'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = '6.47'
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D21").Value = '71.17'
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '11.37'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").Value = '242.91'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("D24").Value = '2.76'
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").Value = '1.96'
$ws.Range("E25").Value = '  +2.86%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '25.27'
$ws.Range("E27").Value = '  +3.11%  '
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("D29").Value = '9.79'
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").Value = '37.15'
$ws.Range("E30").Value = '  +5.31%  '
$ws.Range("D31").Value = '157.62'
$ws.Range("E31").Value = '  +1.35%  '
$ws.Range("D32").Value = '5.52'
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '0.0768'
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = '17.56'
$ws.Range("E36").Value = '  +2.68%  '
$ws.Range("D37").Value = '1.88'
$ws.Range("D38").Value = '2.93'
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").Value = '  +2.12%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = '2.54'
$ws.Range("E42").Value = '  +8.15%  '
$ws.Range("D43").Value = '2.004.94'
$ws.Range("E43").Value = '  +3.35%  '
$ws.Range("D44").Value = '19.43'
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").Value = '2.99'
$ws.Range("E46").Value = '  +3.71%  '
$ws.Range("D47").Value = '9.53'
$ws.Range("E47").Value = '  +5.54%  '
$ws.Range("D48").Value = '2.719.69'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").Value = '76.85'
$ws.Range("E49").Value = '  +7.78%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '98.27'
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("D51").Value = '67.66'
$ws.Range("E51").Value = '  +1.64%  '
